$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 33; $r++) {
    foreach ($col in @("G", "H")) {
        $cell = $ws.Range("$col$r")
        $val = $cell.Value()
        if ($val -ne $null -and $val -ne "N/A") {
            $parts = $val -split ":"
            for ($i = 0; $i -lt $parts.Length; $i++) {
                if ($parts[$i] -like "0x*") {
                    $parts[$i] = "0x" + $parts[$i].Substring(2).ToUpper()
                }
            }
            $cell.Value = [string]::Join(":", $parts)
        }
    }
}
